$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.105721354484558
$ws.Range("B1").Value = 1.500784158706665
$ws.Range("C1").Value = 9.131790161132812
$ws.Range("D1").Value = 2.39056658744812
$ws.Range("E1").Value = 1.283353090286255
